$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "IX Convocatoria de Estímulos a la Excelencia"
$ws.Range("B2").Value = "Dic. 2022"
$ws.Range("C2").Value = "Universidad El Bosque"
$ws.Range("D2").Value = "Bogotá, Colombia"
$ws.Range("E2").Value = "COP`$10.000.000"
$ws.Range("E2").Style = $ws.Range("E4").Style

$ws.Range("D3").Value = "Cambridge, MA, EE.UU."
$ws.Range("E3").Value = "Por ‘tratar de cuantificar la relación entre la desigualdad de ingresos nacionales en diferentes países y la cantidad promedio de besos boca a boca’ (Watkins, et al., 2019)"

$ws.Range("E3").Select()
$excel.ActiveWindow.ScrollColumn = 2
Write-Host "ScrollColumn now:" $excel.ActiveWindow.ScrollColumn
Write-Host "ScrollRow now:" $excel.ActiveWindow.ScrollRow

Write-Host "Inserted row. Dim now:" $ws.UsedRange.Address()
